# Update prediction / stats data on Sheet1 of the Eintracht Frankfurt
# stats workbook (player rows 4-25) to match refreshed sofascore figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("K4").Value = 1.2873

# Row 6
$ws.Range("K6").Value = 0.1254

# Row 9
$ws.Range("E9").Value  = 6.5125
$ws.Range("AB9").Value = 482
$ws.Range("AF9").Value = 84.383561643836
$ws.Range("AG9").Value = 365
$ws.Range("BG9").Value = 106
$ws.Range("BQ9").Value = 104.2
$ws.Range("BT9").Value = 57
$ws.Range("DF9").Value = 170

# Row 10
$ws.Range("K10").Value  = 2.5226
$ws.Range("AA10").Value = 4.07322497
$ws.Range("AR10").Value = 62

# Row 12
$ws.Range("E12").Value  = 6.6866666666667
$ws.Range("K12").Value  = 0.8017
$ws.Range("AF12").Value = 82.135922330097
$ws.Range("AG12").Value = 515
$ws.Range("AO12").Value = 12
$ws.Range("BQ12").Value = 100.3
$ws.Range("BT12").Value = 92
$ws.Range("DF12").Value = 242

# Row 13
$ws.Range("AA13").Value = 2.13849284

# Row 14
$ws.Range("K14").Value = 3.3458

# Row 15
$ws.Range("K15").Value = 0.715

# Row 19
$ws.Range("AB19").Value = 1745
$ws.Range("BG19").Value = 157
$ws.Range("DA19").Value = 12
$ws.Range("DB19").Value = 48

# Row 20
$ws.Range("K20").Value  = 2.0857
$ws.Range("AA20").Value = 2.12604843
$ws.Range("AE20").Value = 541
$ws.Range("AF20").Value = 82.721712538226
$ws.Range("AI20").Value = 237
$ws.Range("AJ20").Value = 116
$ws.Range("AR20").Value = 70
$ws.Range("BG20").Value = 210
$ws.Range("BT20").Value = 113

# Row 21
$ws.Range("AA21").Value = 1.13840623

# Row 25
$ws.Range("E25").Value  = 6.4545454545455
$ws.Range("AB25").Value = 521
$ws.Range("BQ25").Value = 71
$ws.Range("CC25").Value = 31
$ws.Range("CL25").Value = 5
$ws.Range("CO25").Value = 10
$ws.Range("DL25").Value = -8.2836
